$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Region-Country-Mapping")

# Insert a new row at 123 (shifts the Kyrgyzstan..Zimbabwe block, plus the
# trailing Kosovo row, down by one), then move the Kosovo entry (which is
# now the very last row of the sheet) up into the freshly inserted row,
# and attach a footnote explaining the missing ISO3/ISO2 codes in column F.
$ws.Rows.Item(123).Insert()

$lastRow = $ws.Cells.Item($ws.Rows.Count, "B").End(-4162).Row

$ws.Range("B123").Value = $ws.Range("B" + $lastRow).Value()
$ws.Range("C123").Value = $ws.Range("C" + $lastRow).Value()
$ws.Range("F123").Value = "Kosovo does not have ISO3/ISO2 codes under ISO 3166-1 because it is not a universally recognized state"

$ws.Rows.Item($lastRow).Delete()

# Move the active tab from Region-Country-Mapping back to Model.
$wb.Worksheets.Item("Model").Select()
